# Auto-generated script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value is a plain decimal number-looking string ---
# These must be forced to Text format, otherwise Excel auto-converts them
# to numeric values (losing trailing zeros / exact text). We set NumberFormat
# to Text ("@"), assign the value, then reset the style back to "Normal" so
# the cell keeps the default (unstyled) appearance exactly like the rest of
# the sheet.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D5").Value = "300.50"
$ws.Range("D6").Value = "97.85"
$ws.Range("D9").Value = "0.513"
$ws.Range("D10").Value = "35.00"
$ws.Range("D11").Value = "0.0794"
$ws.Range("D12").Value = "7.19"
$ws.Range("D17").Value = "14.21"
$ws.Range("D19").Value = "13.09"
$ws.Range("D21").Value = "6.25"
$ws.Range("D22").Value = "67.22"
$ws.Range("D23").Value = "243.46"
$ws.Range("D25").Value = "1.00"
$ws.Range("D26").Value = "1.94"
$ws.Range("D27").Value = "38.55"
$ws.Range("D29").Value = "9.82"
$ws.Range("D33").Value = "5.56"
$ws.Range("D34").Value = "148.54"
$ws.Range("D35").Value = "0.0778"
$ws.Range("D36").Value = "1.98"
$ws.Range("D37").Value = "0.114"
$ws.Range("D39").Value = "15.15"
$ws.Range("D42").Value = "3.29"
$ws.Range("D44").Value = "0.998"
$ws.Range("D45").Value = "91.60"
$ws.Range("D46").Value = "1.76"
$ws.Range("D48").Value = "103.29"
$ws.Range("D50").Value = "0.189"

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"

# --- Remaining cells: safe to assign directly since Excel will not ---
# --- reinterpret them as numbers (they contain letters, extra dots, ---
# --- percent signs/spaces, etc.) ---
$ws.Range("D2").Value = "45.650.00"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").Value = "2.409.41"
$ws.Range("E3").Value = "  +5.12%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("E6").Value = "  -2.11%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("E10").Value = "  -3.52%  "
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("E12").Value = "  -2.72%  "
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "2.761.05"
$ws.Range("E14").Value = "  +4.43%  "
$ws.Range("D15").Value = "2.414.21"
$ws.Range("E15").Value = "  +5.27%  "
$ws.Range("E16").Value = "  +4.35%  "
$ws.Range("E17").Value = "  +3.04%  "
$ws.Range("D18").Value = "45.599.76"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("E23").Value = "  -2.04%  "
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("E27").Value = "  -9.27%  "
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("E30").Value = "  +16.87%  "
$ws.Range("E31").Value = "  +6.63%  "
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("E35").Value = "  -1.63%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E36").Value = "  +12.07%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("E39").Value = "  -4.84%  "
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").Value = "1.950.87"
$ws.Range("E43").Value = "  +7.13%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("E45").Value = "  +4.15%  "
$ws.Range("E46").Value = "  -11.00%  "
$ws.Range("E47").Value = "  +10.24%  "
$ws.Range("E48").Value = "  +7.88%  "
$ws.Range("E49").Value = "  +14.79%  "
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("D51").Value = "2.645.25"
$ws.Range("E51").Value = "  +4.96%  "
